$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they stay text like the rest of the column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '62.039.36'
$ws.Range("E2").Value = '  -0.32%  '

# Row 3
$ws.Range("D3").Value = '2.418.21'
$ws.Range("E3").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '563.77'
$ws.Range("E5").Value = '  +1.37%  '

# Row 6
$ws.Range("D6").Value = '142.78'
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("E8").Value = '  -0.25%  '

# Row 9
$ws.Range("E9").Value = '  +0.33%  '

# Row 10
$ws.Range("E10").Value = '  -1.84%  '

# Row 11
$ws.Range("E11").Value = '  -3.05%  '

# Row 12
$ws.Range("E12").Value = '  -0.31%  '

# Row 13
$ws.Range("D13").Value = '25.62'
$ws.Range("E13").Value = '  -2.26%  '

# Row 14
$ws.Range("E14").Value = '  +0.15%  '

# Row 15
$ws.Range("D15").Value = '2.854.96'
$ws.Range("E15").Value = '  +0.30%  '

# Row 16
$ws.Range("D16").Value = '61.915.05'
$ws.Range("E16").Value = '  -0.29%  '

# Row 17
$ws.Range("D17").Value = '2.420.78'
$ws.Range("E17").Value = '  +0.14%  '

# Row 18
$ws.Range("D18").Value = '11.24'
$ws.Range("E18").Value = '  +1.53%  '

# Row 19
$ws.Range("D19").Value = '321.58'
$ws.Range("E19").Value = '  -0.71%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '6.83'
$ws.Range("E20").Value = '  +1.80%  '

# Row 21
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '4.14'
$ws.Range("E21").Value = '  -1.35%  '

# Row 22
$ws.Range("E22").Value = '  -0.14%  '

# Row 23
$ws.Range("D23").Value = '66.09'
$ws.Range("E23").Value = '  +1.78%  '

# Row 24
$ws.Range("E24").Value = '  -0.46%  '

# Row 25
$ws.Range("E25").Value = '  -5.03%  '

# Row 26
$ws.Range("D26").Value = '564.76'
$ws.Range("E26").Value = '  -2.21%  '

# Row 27
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.529.79'
$ws.Range("E27").Value = '  -0.48%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.02%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0934'
$ws.Range("E29").Value = '  -0.01%  '

# Row 30
$ws.Range("E30").Value = '  -2.20%  '

# Row 31
$ws.Range("E31").Value = '  -3.70%  '

# Row 32
$ws.Range("E32").Value = '  -0.24%  '

# Row 33
$ws.Range("E33").Value = '  +1.36%  '

# Row 34
$ws.Range("E34").Value = '  -2.03%  '

# Row 35
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.01%  '

# Row 36
$ws.Range("E36").Value = '  -0.94%  '

# Row 37
$ws.Range("E37").Value = '  -3.44%  '

# Row 38
$ws.Range("D38").Value = '153.06'
$ws.Range("E38").Value = '  +3.15%  '

# Row 39
$ws.Range("E39").Value = '  -1.05%  '

# Row 40
$ws.Range("D40").Value = '18.55'
$ws.Range("E40").Value = '  -0.97%  '

# Row 41
$ws.Range("E41").Value = '  -3.71%  '

# Row 42
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("D43").Value = '148.91'
$ws.Range("E43").Value = '  -1.24%  '

# Row 44
$ws.Range("E44").Value = '  -2.35%  '

# Row 45
$ws.Range("D45").Value = '3.62'
$ws.Range("E45").Value = '  -0.41%  '

# Row 46
$ws.Range("E46").Value = '  -2.23%  '

# Row 47
$ws.Range("E47").Value = '  -2.34%  '

# Row 48
$ws.Range("E48").Value = '  +0.75%  '

# Row 49
$ws.Range("E49").Value = '  +0.23%  '

# Row 50
$ws.Range("E50").Value = '  -0.85%  '

# Row 51
$ws.Range("E51").Value = '  +0.64%  '
